$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.567.96"
$ws.Range("E2").Value = "'  -1.43%  "
$ws.Range("D3").Value = "'2.033.91"
$ws.Range("E3").Value = "'  +1.49%  "
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'234.26"
$ws.Range("E5").Value = "'  -8.98%  "
$ws.Range("D6").Value = "'0.595"
$ws.Range("E6").Value = "'  -3.39%  "
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("D8").Value = "'55.19"
$ws.Range("E8").Value = "'  -1.24%  "
$ws.Range("D9").Value = "'0.372"
$ws.Range("E9").Value = "'  -1.49%  "
$ws.Range("D10").Value = "'57.51"
$ws.Range("E10").Value = "'  +2.70%  "
$ws.Range("D11").Value = "'0.0753"
$ws.Range("E11").Value = "'  -1.74%  "
$ws.Range("E12").Value = "'  -1.60%  "
$ws.Range("D13").Value = "'2.329.88"
$ws.Range("E13").Value = "'  +1.09%  "
$ws.Range("D14").Value = "'14.33"
$ws.Range("E14").Value = "'  +0.67%  "
$ws.Range("D15").Value = "'20.37"
$ws.Range("E15").Value = "'  -4.70%  "
$ws.Range("E16").Value = "'  -3.71%  "
$ws.Range("D17").Value = "'5.10"
$ws.Range("E17").Value = "'  -1.48%  "
$ws.Range("D18").Value = "'2.035.17"
$ws.Range("E18").Value = "'  +0.59%  "
$ws.Range("D19").Value = "'36.706.37"
$ws.Range("E19").Value = "'  -0.62%  "
$ws.Range("D20").Value = "'67.81"
$ws.Range("E20").Value = "'  -4.26%  "
$ws.Range("B21").Value = "'Uniswap"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.48"
$ws.Range("E21").Value = "'  +8.05%  "
$ws.Range("B22").Value = "'ShibaInu"
$ws.Range("C22").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "'0.0₃0801"
$ws.Range("E22").Value = "'  -3.36%  "
$ws.Range("D23").Value = "'221.37"
$ws.Range("E23").Value = "'  -5.18%  "
$ws.Range("E24").Value = "'  +0.08%  "
$ws.Range("E25").Value = "'  +1.10%  "
$ws.Range("E26").Value = "'  -5.67%  "
$ws.Range("D27").Value = "'163.21"
$ws.Range("E27").Value = "'  -0.79%  "
$ws.Range("E28").Value = "'  +8.23%  "
$ws.Range("D29").Value = "'8.65"
$ws.Range("E29").Value = "'  -2.83%  "
$ws.Range("D30").Value = "'19.02"
$ws.Range("E30").Value = "'  -2.25%  "
$ws.Range("D31").Value = "'1.36"
$ws.Range("E31").Value = "'  +1.55%  "
$ws.Range("E32").Value = "'  -1.63%  "
$ws.Range("D33").Value = "'4.39"
$ws.Range("E33").Value = "'  -3.96%  "
$ws.Range("D34").Value = "'0.0604"
$ws.Range("E34").Value = "'  -5.30%  "
$ws.Range("E35").Value = "'  +5.65%  "
$ws.Range("D36").Value = "'4.27"
$ws.Range("E36").Value = "'  -2.90%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "'  -0.30%  "
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "'  -4.04%  "
$ws.Range("E39").Value = "'  -2.90%  "
$ws.Range("D40").Value = "'5.79"
$ws.Range("E40").Value = "'  +4.09%  "
$ws.Range("E41").Value = "'  -4.28%  "
$ws.Range("D42").Value = "'0.0948"
$ws.Range("E42").Value = "'  +3.30%  "
$ws.Range("D43").Value = "'1.460.83"
$ws.Range("E43").Value = "'  +1.59%  "
$ws.Range("B44").Value = "'FTXToken"
$ws.Range("C44").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.19"
$ws.Range("E44").Value = "'  +41.23%  "
$ws.Range("B45").Value = "'VeChain"
$ws.Range("C45").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0205"
$ws.Range("E45").Value = "'  -2.34%  "
$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'90.99"
$ws.Range("E46").Value = "'  +2.49%  "
$ws.Range("B47").Value = "'TrustWalletToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "'1.11"
$ws.Range("E47").Value = "'  -5.26%  "
$ws.Range("D48").Value = "'15.63"
$ws.Range("E48").Value = "'  +0.85%  "
$ws.Range("E49").Value = "'  -0.68%  "
$ws.Range("E50").Value = "'  -1.40%  "
$ws.Range("D51").Value = "'6.89"
$ws.Range("E51").Value = "'  -0.39%  "
